# Framework Setup - App elements, Steps and Tests
#
# Replaces the generic Sheet1/Sheet2 smoke-test data with the real
# OrangeHRM TestNG framework TestData workbook: 5 named sheets, each
# holding the login credentials / header text used by the automated
# test cases.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Build the 5 target sheets (appended after the existing Sheet2 so
#    Add()'s internal sheetId counter lands on 3,4,5,6,7 - matching the
#    sequence the original author ended up with), then drop the two
#    stock sheets and reorder so the id-7 sheet sits in slot 3.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sh1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$sh1.Name = "verifyLoginPageLogoAndHeader"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sh2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$sh2.Name = "verifyLoginWithValidCredentials"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sh4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$sh4.Name = "VerifyProfileImageInHomePage"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sh5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$sh5.Name = "VerifyHomePageMenus"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sh3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$sh3.Name = "verifyInvalidLogin"

$wb.Worksheets.Item("Sheet1").Delete() | Out-Null
$wb.Worksheets.Item("Sheet2").Delete() | Out-Null

# Put verifyInvalidLogin (3rd tab) back between verifyLoginWithValidCredentials
# and VerifyProfileImageInHomePage.
$wb.Worksheets.Item("verifyInvalidLogin").Move($wb.Worksheets.Item("VerifyProfileImageInHomePage"))

# ---------------------------------------------------------------------
# 2) verifyLoginPageLogoAndHeader - the page header/logo text checked by
#    the first test case.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("verifyLoginPageLogoAndHeader")
$ws.Range("A1").Value = "LoginHeader"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A2").Value = "Login"
$ws.Columns.Item(1).ColumnWidth = 11.33
$ws.Range("J32").Select() | Out-Null

# ---------------------------------------------------------------------
# 3) verifyLoginWithValidCredentials - Admin / admin123 credential pair.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("verifyLoginWithValidCredentials")
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "admin123"
$ws.Columns.Item(1).ColumnWidth = 9.33
$ws.Columns.Item(2).ColumnWidth = 10
$ws.Rows.Item(3).Select() | Out-Null

# ---------------------------------------------------------------------
# 4) verifyInvalidLogin - Bharath / bharath123 (bad) credential pair.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("verifyInvalidLogin")
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A2").Value = "Bharath"
$ws.Range("B2").Value = "bharath123"
$ws.Columns.Item(1).ColumnWidth = 9.33
$ws.Columns.Item(2).ColumnWidth = 10
$ws.Activate()
$ws.Range("F21").Select() | Out-Null

# ---------------------------------------------------------------------
# 5) VerifyProfileImageInHomePage - re-uses the Admin / admin123 pair.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("VerifyProfileImageInHomePage")
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "admin123"
$ws.Range("G37").Select() | Out-Null

# ---------------------------------------------------------------------
# 6) VerifyHomePageMenus - re-uses the Admin / admin123 pair.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("VerifyHomePageMenus")
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "admin123"
$ws.Range("Q28").Select() | Out-Null

# ---------------------------------------------------------------------
# 7) Make verifyInvalidLogin the active/visible tab, matching the saved
#    workbook view (activeTab points at the 3rd sheet).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("verifyInvalidLogin").Activate()

foreach ($ws in $wb.Worksheets) { Write-Output $ws.Name }
